# Append four blank paragraphs to the very end of the document body,
# right before the final section break (sectPr). Three of them contain
# a single empty run, and the last one is completely empty, matching:
#
#   <w:p><w:r/></w:p>
#   <w:p><w:r/></w:p>
#   <w:p><w:r/></w:p>
#   <w:p/>
#
# Inserting raw WordprocessingML via Range.InsertXML (rather than the
# higher-level InsertParagraphAfter/TypeParagraph calls) avoids pulling
# in any inherited run/paragraph formatting, so the new paragraphs stay
# perfectly empty just like the target markup.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$blankParagraphsXml = "<w:p $wNs><w:r/></w:p>" + `
                       "<w:p $wNs><w:r/></w:p>" + `
                       "<w:p $wNs><w:r/></w:p>" + `
                       "<w:p $wNs/>"

$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)
[void]$insertionPoint.InsertXML($blankParagraphsXml)
